# Update Name of Algo
# Applies updated numeric results for the RandomForest imputation output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = -8.124499999999994
$ws.Range("D6").Value = -8.729999999999995
$ws.Range("C7").Value = -12.2878
$ws.Range("A10").Value = -20.48899999999997
$ws.Range("A12").Value = -22.52250000000004
$ws.Range("B13").Value = 5.921699999999996
$ws.Range("A18").Value = -22.45250000000003
$ws.Range("C20").Value = -14.61409999999999
